$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data cells for rows 2-5 (new sensor readings, dates shifted +279 days) ---
$row2 = @(45060.50694444445,4.928,5.344,0,5.978,10.253,2.511,8.567,4.041,1.748,4.53,5.048,5.182,0.791,3.479,3.593,1.503,1.042,0.417,43.098,7.946,4.36,6.645,2.748,0.484,2.765,1.536,3.068,2.75,5.612,0,5.035,2.082,3.31)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}
$row3 = @(45060.51388888889,11.375,9.19,0.185,22.54,21.065,8.185,29.461,12.777,5.845,9.587,10.127,10.653,2.666,8.689,11.773,6.449,0.681,0.555,123.872,23.487,8.48,16.702,8.406,1.147,14.165,6.43,6.467,7.414,10.904,0,25.477,4.777,9.665)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}
$row4 = @(45060.52083333334,14.632,11.405,0.331,30.397,26.708,10.978,42.991,17.053,7.771,12.131,12.886,13.581,3.564,11.296,15.744,8.856,0.53,0.599,164.1,31.13,10.711,21.652,11.131,1.485,20.596,8.777,8.239,9.625,13.72,0,38.335,6.108,12.803)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}
$row5 = @(45060.52777777778,19.38,14.87,0.52,41.1,35.15,14.84,58.2,22.97,10.42,15.93,17,17.94,4.79,15.04,21.22,12.09,0.47,0.71,221.4,41.74,14.09,28.72,14.91,1.98,27.99,11.94,10.88,12.77,18,0,52.26,8.05,17.19)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5[$i]
}

# --- Delete row 6 (dataset trimmed to 4 data rows) ---
$ws.Rows.Item(6).Delete()

# --- Column width changes (ColumnWidth = XML stored width - 0.83) ---
$ws.Range("B:B").ColumnWidth = 7.17
$ws.Range("C:C").ColumnWidth = 7.17
$ws.Range("G:G").ColumnWidth = 7.17
$ws.Range("I:I").ColumnWidth = 7.17
$ws.Range("K:K").ColumnWidth = 7.17
$ws.Range("L:L").ColumnWidth = 7.17
$ws.Range("M:M").ColumnWidth = 7.17
$ws.Range("O:O").ColumnWidth = 7.17
$ws.Range("P:P").ColumnWidth = 7.17
$ws.Range("T:T").ColumnWidth = 8.17
$ws.Range("V:V").ColumnWidth = 7.17
$ws.Range("X:X").ColumnWidth = 7.17
$ws.Range("Z:Z").ColumnWidth = 7.17
$ws.Range("AD:AD").ColumnWidth = 7.17
$ws.Range("AE:AE").ColumnWidth = 4.17
$ws.Range("AH:AH").ColumnWidth = 7.17
